$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C1 value (rounding correction for Misc.significantNumbers method)
$ws.Range("C1").Value = 123.21312

# Add new imported data in row 2 (excel data import fix)
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 123
$ws.Range("F2").Value = 123

# Remove row 7 (C7 = 5) entirely
$ws.Rows("7:7").Delete()

# Select cell F3 as the active cell
$ws.Range("F3").Select()
